# Auto-generated Excel COM-interop script to apply profit/cost recalculation updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 78
$ws.Range("I2").Value = 78
$ws.Range("K2").Value = 78
$ws.Range("M2").Value = 35

$ws.Range("H18").Value = 4350
$ws.Range("I18").Value = 3700
$ws.Range("K18").Value = 3700
$ws.Range("M18").Value = -3416

$ws.Range("H19").Value = 898.25
$ws.Range("I19").Value = 898
$ws.Range("J19").Value = 899
$ws.Range("K19").Value = 898
$ws.Range("L19").Value = 899
$ws.Range("M19").Value = -723
$ws.Range("N19").Value = -1249

$ws.Range("H29").Value = 181.625
$ws.Range("I29").Value = 181.625
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 544.875
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -263.875
$ws.Range("N29").ClearContents()

$ws.Range("H38").Value = 587
$ws.Range("I38").Value = 587
$ws.Range("K38").Value = 1761
$ws.Range("M38").Value = -1389

$ws.Range("H51").Value = 4000
$ws.Range("I51").Value = 4000
$ws.Range("K51").Value = 4000
$ws.Range("M51").Value = -3516

$ws.Range("H80").Value = 567.1667
$ws.Range("I80").Value = 650
$ws.Range("J80").Value = 525.75
$ws.Range("K80").Value = 1950
$ws.Range("L80").Value = 1577.25
$ws.Range("M80").Value = -952
$ws.Range("N80").Value = -3573.25

$ws.Range("H83").Value = 567.1667
$ws.Range("I83").Value = 650
$ws.Range("J83").Value = 525.75
$ws.Range("K83").Value = 5850
$ws.Range("L83").Value = 4731.75
$ws.Range("M83").Value = -858
$ws.Range("N83").Value = -14715.75

$ws.Range("H98").Value = 17211.5
$ws.Range("I98").Value = 17211.5
$ws.Range("K98").Value = 17211.5
$ws.Range("M98").Value = -15713.5

$ws.Range("H100").Value = 1891.375
$ws.Range("I100").Value = 1018.8571
$ws.Range("K100").Value = 1018.8571
$ws.Range("M100").Value = -477.8570999999999

$ws.Range("H111").Value = 3277
$ws.Range("I111").Value = 902.1667
$ws.Range("J111").Value = 5312.5713
$ws.Range("K111").Value = 2706.5001
$ws.Range("L111").Value = 15937.7139
$ws.Range("M111").Value = 360.4998999999998
$ws.Range("N111").Value = -22071.7139

$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 2254
$ws.Range("N113").Value = -7708

$ws.Range("H122").Value = 17211.5
$ws.Range("I122").Value = 17211.5
$ws.Range("K122").Value = 51634.5
$ws.Range("M122").Value = -49184.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 6332.8335
$ws.Range("J110").Value = 6499.5
$ws.Range("L110").Value = 6499.5
$ws.Range("N110").Value = -10589.5

$ws.Range("H111").Value = 65000
$ws.Range("J111").Value = 65000
$ws.Range("L111").Value = 65000
$ws.Range("N111").Value = -73180

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2206
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H107").Value = 16588.375
$ws.Range("I107").Value = 7701.3335
$ws.Range("J107").Value = 43249.5
$ws.Range("K107").Value = 7701.3335
$ws.Range("L107").Value = 43249.5
$ws.Range("M107").Value = -5781.3335
$ws.Range("N107").Value = -47089.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 869.5
$ws.Range("I16").Value = 869.5
$ws.Range("K16").Value = 869.5
$ws.Range("M16").Value = -582.5

$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498

$ws.Range("H107").Value = 388.69232
$ws.Range("I107").Value = 356.5
$ws.Range("K107").Value = 356.5
$ws.Range("M107").Value = 1563.5

$ws.Range("H113").Value = 869.5
$ws.Range("I113").Value = 869.5
$ws.Range("K113").Value = 869.5
$ws.Range("M113").Value = 1300.5

$ws.Range("H122").Value = 1706.8572
$ws.Range("J122").Value = 1499.5
$ws.Range("L122").Value = 4498.5
$ws.Range("N122").Value = -9398.5

$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 17
$ws.Range("I12").Value = 8.285714
$ws.Range("J12").Value = 24.625
$ws.Range("K12").Value = 24.857142
$ws.Range("L12").Value = 73.875
$ws.Range("M12").Value = 148.142858
$ws.Range("N12").Value = -419.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 59944.555
$ws.Range("J69").Value = 59944.555
$ws.Range("L69").Value = 59944.555
$ws.Range("N69").Value = -61442.555

$ws.Range("H72").Value = 59944.555
$ws.Range("J72").Value = 59944.555
$ws.Range("L72").Value = 179833.665
$ws.Range("N72").Value = -187321.665

$ws.Range("H102").Value = 4240
$ws.Range("I102").Value = 5075
$ws.Range("K102").Value = 5075
$ws.Range("M102").Value = -3453

$ws.Range("H122").Value = 24337.5
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 47425
$ws.Range("K122").Value = 3750
$ws.Range("L122").Value = 142275
$ws.Range("M122").Value = -1300
$ws.Range("N122").Value = -147175

$ws.Range("H126").Value = 4247
$ws.Range("I126").Value = 3329.3333
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 9987.999899999999
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -7517.999899999999
$ws.Range("N126").Value = -25940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2783.5
$ws.Range("I61").Value = 3179.2
$ws.Range("K61").Value = 3179.2
$ws.Range("M61").Value = -2977.2

$ws.Range("H113").Value = 2783.5
$ws.Range("I113").Value = 3179.2
$ws.Range("K113").Value = 3179.2
$ws.Range("M113").Value = -1009.2

$ws.Range("H136").Value = 17999.75
$ws.Range("I136").Value = 18500
$ws.Range("K136").Value = 55500
$ws.Range("M136").Value = -52950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 919.94116
$ws.Range("I107").Value = 568.25
$ws.Range("K107").Value = 1704.75
$ws.Range("M107").Value = 215.25

$ws.Range("H122").Value = 224424.56
$ws.Range("I122").Value = 251602.62
$ws.Range("K122").Value = 754807.86
$ws.Range("M122").Value = -752357.86
